$d = $word.ActiveDocument
$d.Content.Find.Execute("97×29=2813", $true, $false, $false, $false, $false, $true, 1, $false, "20×22=440", 2) | Out-Null
$d.Content.Find.Execute("94×73=6862", $true, $false, $false, $false, $false, $true, 1, $false, "49×54=2646", 2) | Out-Null
$d.Content.Find.Execute("48×94=4512", $true, $false, $false, $false, $false, $true, 1, $false, "86×93=7998", 2) | Out-Null
$d.Content.Find.Execute("81×12=972", $true, $false, $false, $false, $false, $true, 1, $false, "42×15=630", 2) | Out-Null
$d.Content.Find.Execute("15×69=1035", $true, $false, $false, $false, $false, $true, 1, $false, "51×98=4998", 2) | Out-Null
$d.Content.Find.Execute("93×56=5208", $true, $false, $false, $false, $false, $true, 1, $false, "74×78=5772", 2) | Out-Null
$d.Content.Find.Execute("31×23=713", $true, $false, $false, $false, $false, $true, 1, $false, "59×42=2478", 2) | Out-Null
$d.Content.Find.Execute("51×51=2601", $true, $false, $false, $false, $false, $true, 1, $false, "86×33=2838", 2) | Out-Null
$d.Content.Find.Execute("95×37=3515", $true, $false, $false, $false, $false, $true, 1, $false, "13×35=455", 2) | Out-Null
$d.Content.Find.Execute("48×35=1680", $true, $false, $false, $false, $false, $true, 1, $false, "29×56=1624", 2) | Out-Null
$d.Content.Find.Execute("93×18=1674", $true, $false, $false, $false, $false, $true, 1, $false, "85×68=5780", 2) | Out-Null
$d.Content.Find.Execute("45×33=1485", $true, $false, $false, $false, $false, $true, 1, $false, "57×35=1995", 2) | Out-Null
$d.Content.Find.Execute("76×66=5016", $true, $false, $false, $false, $false, $true, 1, $false, "58×42=2436", 2) | Out-Null
$d.Content.Find.Execute("88×84=7392", $true, $false, $false, $false, $false, $true, 1, $false, "62×83=5146", 2) | Out-Null
$d.Content.Find.Execute("66×47=3102", $true, $false, $false, $false, $false, $true, 1, $false, "74×87=6438", 2) | Out-Null
$d.Content.Find.Execute("32×48=1536", $true, $false, $false, $false, $false, $true, 1, $false, "21×62=1302", 2) | Out-Null
$d.Content.Find.Execute("68×93=6324", $true, $false, $false, $false, $false, $true, 1, $false, "35×50=1750", 2) | Out-Null
$d.Content.Find.Execute("40×14=560", $true, $false, $false, $false, $false, $true, 1, $false, "73×24=1752", 2) | Out-Null
$d.Content.Find.Execute("27×11=297", $true, $false, $false, $false, $false, $true, 1, $false, "83×77=6391", 2) | Out-Null
$d.Content.Find.Execute("34×29=986", $true, $false, $false, $false, $false, $true, 1, $false, "33×27=891", 2) | Out-Null
$d.Content.Find.Execute("67×12=804", $true, $false, $false, $false, $false, $true, 1, $false, "60×84=5040", 2) | Out-Null
$d.Content.Find.Execute("61×58=3538", $true, $false, $false, $false, $false, $true, 1, $false, "20×51=1020", 2) | Out-Null
$d.Content.Find.Execute("98×47=4606", $true, $false, $false, $false, $false, $true, 1, $false, "48×24=1152", 2) | Out-Null
$d.Content.Find.Execute("34×95=3230", $true, $false, $false, $false, $false, $true, 1, $false, "75×50=3750", 2) | Out-Null
$d.Content.Find.Execute("30×53=1590", $true, $false, $false, $false, $false, $true, 1, $false, "25×19=475", 2) | Out-Null
$d.Content.Find.Execute("29×79=2291", $true, $false, $false, $false, $false, $true, 1, $false, "98×23=2254", 2) | Out-Null
$d.Content.Find.Execute("60×17=1020", $true, $false, $false, $false, $false, $true, 1, $false, "71×98=6958", 2) | Out-Null
$d.Content.Find.Execute("67×70=4690", $true, $false, $false, $false, $false, $true, 1, $false, "78×12=936", 2) | Out-Null
$d.Content.Find.Execute("68×57=3876", $true, $false, $false, $false, $false, $true, 1, $false, "66×16=1056", 2) | Out-Null
$d.Content.Find.Execute("41×72=2952", $true, $false, $false, $false, $false, $true, 1, $false, "79×26=2054", 2) | Out-Null
$d.Content.Find.Execute("55×69=3795", $true, $false, $false, $false, $false, $true, 1, $false, "75×47=3525", 2) | Out-Null
$d.Content.Find.Execute("82×76=6232", $true, $false, $false, $false, $false, $true, 1, $false, "55×36=1980", 2) | Out-Null
$d.Content.Find.Execute("58×61=3538", $true, $false, $false, $false, $false, $true, 1, $false, "54×31=1674", 2) | Out-Null
$d.Content.Find.Execute("23×98=2254", $true, $false, $false, $false, $false, $true, 1, $false, "92×79=7268", 2) | Out-Null
$d.Content.Find.Execute("82×18=1476", $true, $false, $false, $false, $false, $true, 1, $false, "32×95=3040", 2) | Out-Null
$d.Content.Find.Execute("68×79=5372", $true, $false, $false, $false, $false, $true, 1, $false, "61×53=3233", 2) | Out-Null
$d.Content.Find.Execute("49×48=2352", $true, $false, $false, $false, $false, $true, 1, $false, "44×91=4004", 2) | Out-Null
$d.Content.Find.Execute("60×76=4560", $true, $false, $false, $false, $false, $true, 1, $false, "66×68=4488", 2) | Out-Null
$d.Content.Find.Execute("24×28=672", $true, $false, $false, $false, $false, $true, 1, $false, "67×83=5561", 2) | Out-Null
$d.Content.Find.Execute("44×22=968", $true, $false, $false, $false, $false, $true, 1, $false, "21×23=483", 2) | Out-Null
$d.Content.Find.Execute("56×76=4256", $true, $false, $false, $false, $false, $true, 1, $false, "72×74=5328", 2) | Out-Null
$d.Content.Find.Execute("10×54=540", $true, $false, $false, $false, $false, $true, 1, $false, "87×71=6177", 2) | Out-Null
$d.Content.Find.Execute("91×83=7553", $true, $false, $false, $false, $false, $true, 1, $false, "94×52=4888", 2) | Out-Null
$d.Content.Find.Execute("54×84=4536", $true, $false, $false, $false, $false, $true, 1, $false, "67×44=2948", 2) | Out-Null
$d.Content.Find.Execute("100×20=2000", $true, $false, $false, $false, $false, $true, 1, $false, "57×92=5244", 2) | Out-Null
$d.Content.Find.Execute("11×65=715", $true, $false, $false, $false, $false, $true, 1, $false, "66×38=2508", 2) | Out-Null
$d.Content.Find.Execute("18×71=1278", $true, $false, $false, $false, $false, $true, 1, $false, "100×12=1200", 2) | Out-Null
$d.Content.Find.Execute("22×15=330", $true, $false, $false, $false, $false, $true, 1, $false, "15×89=1335", 2) | Out-Null
$d.Content.Find.Execute("18×79=1422", $true, $false, $false, $false, $false, $true, 1, $false, "41×10=410", 2) | Out-Null
$d.Content.Find.Execute("25×17=425", $true, $false, $false, $false, $false, $true, 1, $false, "45×43=1935", 2) | Out-Null
$d.Content.Find.Execute("79×51=4029", $true, $false, $false, $false, $false, $true, 1, $false, "58×56=3248", 2) | Out-Null
$d.Content.Find.Execute("58×18=1044", $true, $false, $false, $false, $false, $true, 1, $false, "99×70=6930", 2) | Out-Null
$d.Content.Find.Execute("100×48=4800", $true, $false, $false, $false, $false, $true, 1, $false, "75×38=2850", 2) | Out-Null
$d.Content.Find.Execute("60×35=2100", $true, $false, $false, $false, $false, $true, 1, $false, "88×96=8448", 2) | Out-Null
$d.Content.Find.Execute("35×18=630", $true, $false, $false, $false, $false, $true, 1, $false, "11×94=1034", 2) | Out-Null
$d.Content.Find.Execute("80×82=6560", $true, $false, $false, $false, $false, $true, 1, $false, "83×50=4150", 2) | Out-Null
$d.Content.Find.Execute("16×43=688", $true, $false, $false, $false, $false, $true, 1, $false, "90×16=1440", 2) | Out-Null
$d.Content.Find.Execute("64×94=6016", $true, $false, $false, $false, $false, $true, 1, $false, "97×19=1843", 2) | Out-Null
$d.Content.Find.Execute("66×54=3564", $true, $false, $false, $false, $false, $true, 1, $false, "19×52=988", 2) | Out-Null
$d.Content.Find.Execute("95×57=5415", $true, $false, $false, $false, $false, $true, 1, $false, "13×98=1274", 2) | Out-Null
$d.Content.Find.Execute("97×42=4074", $true, $false, $false, $false, $false, $true, 1, $false, "60×52=3120", 2) | Out-Null
$d.Content.Find.Execute("79×11=869", $true, $false, $false, $false, $false, $true, 1, $false, "20×97=1940", 2) | Out-Null
$d.Content.Find.Execute("16×75=1200", $true, $false, $false, $false, $false, $true, 1, $false, "40×28=1120", 2) | Out-Null
$d.Content.Find.Execute("97×41=3977", $true, $false, $false, $false, $false, $true, 1, $false, "73×71=5183", 2) | Out-Null
$d.Content.Find.Execute("43×65=2795", $true, $false, $false, $false, $false, $true, 1, $false, "46×22=1012", 2) | Out-Null
$d.Content.Find.Execute("25×60=1500", $true, $false, $false, $false, $false, $true, 1, $false, "98×93=9114", 2) | Out-Null
$d.Content.Find.Execute("45×62=2790", $true, $false, $false, $false, $false, $true, 1, $false, "58×22=1276", 2) | Out-Null
$d.Content.Find.Execute("96×36=3456", $true, $false, $false, $false, $false, $true, 1, $false, "58×93=5394", 2) | Out-Null
$d.Content.Find.Execute("89×19=1691", $true, $false, $false, $false, $false, $true, 1, $false, "16×84=1344", 2) | Out-Null
$d.Content.Find.Execute("79×36=2844", $true, $false, $false, $false, $false, $true, 1, $false, "73×90=6570", 2) | Out-Null
$d.Content.Find.Execute("33×32=1056", $true, $false, $false, $false, $false, $true, 1, $false, "55×30=1650", 2) | Out-Null
$d.Content.Find.Execute("85×98=8330", $true, $false, $false, $false, $false, $true, 1, $false, "42×63=2646", 2) | Out-Null
$d.Content.Find.Execute("25×33=825", $true, $false, $false, $false, $false, $true, 1, $false, "11×29=319", 2) | Out-Null
$d.Content.Find.Execute("61×79=4819", $true, $false, $false, $false, $false, $true, 1, $false, "69×30=2070", 2) | Out-Null
$d.Content.Find.Execute("90×11=990", $true, $false, $false, $false, $false, $true, 1, $false, "81×39=3159", 2) | Out-Null
$d.Content.Find.Execute("52×74=3848", $true, $false, $false, $false, $false, $true, 1, $false, "63×50=3150", 2) | Out-Null
$d.Content.Find.Execute("72×76=5472", $true, $false, $false, $false, $false, $true, 1, $false, "18×100=1800", 2) | Out-Null
$d.Content.Find.Execute("25×82=2050", $true, $false, $false, $false, $false, $true, 1, $false, "100×32=3200", 2) | Out-Null
$d.Content.Find.Execute("70×60=4200", $true, $false, $false, $false, $false, $true, 1, $false, "85×15=1275", 2) | Out-Null
$d.Content.Find.Execute("71×34=2414", $true, $false, $false, $false, $false, $true, 1, $false, "85×41=3485", 2) | Out-Null
$d.Content.Find.Execute("22×86=1892", $true, $false, $false, $false, $false, $true, 1, $false, "32×20=640", 2) | Out-Null
$d.Content.Find.Execute("70×40=2800", $true, $false, $false, $false, $false, $true, 1, $false, "97×63=6111", 2) | Out-Null
$d.Content.Find.Execute("62×13=806", $true, $false, $false, $false, $false, $true, 1, $false, "55×62=3410", 2) | Out-Null
$d.Content.Find.Execute("52×68=3536", $true, $false, $false, $false, $false, $true, 1, $false, "51×12=612", 2) | Out-Null
$d.Content.Find.Execute("89×17=1513", $true, $false, $false, $false, $false, $true, 1, $false, "11×64=704", 2) | Out-Null
$d.Content.Find.Execute("61×51=3111", $true, $false, $false, $false, $false, $true, 1, $false, "95×55=5225", 2) | Out-Null
$d.Content.Find.Execute("92×62=5704", $true, $false, $false, $false, $false, $true, 1, $false, "58×83=4814", 2) | Out-Null
$d.Content.Find.Execute("27×73=1971", $true, $false, $false, $false, $false, $true, 1, $false, "51×20=1020", 2) | Out-Null
$d.Content.Find.Execute("99×57=5643", $true, $false, $false, $false, $false, $true, 1, $false, "92×32=2944", 2) | Out-Null
$d.Content.Find.Execute("100×21=2100", $true, $false, $false, $false, $false, $true, 1, $false, "25×21=525", 2) | Out-Null
$d.Content.Find.Execute("83×37=3071", $true, $false, $false, $false, $false, $true, 1, $false, "82×29=2378", 2) | Out-Null
$d.Content.Find.Execute("43×28=1204", $true, $false, $false, $false, $false, $true, 1, $false, "39×39=1521", 2) | Out-Null
$d.Content.Find.Execute("91×53=4823", $true, $false, $false, $false, $false, $true, 1, $false, "48×53=2544", 2) | Out-Null
$d.Content.Find.Execute("66×43=2838", $true, $false, $false, $false, $false, $true, 1, $false, "96×75=7200", 2) | Out-Null
$d.Content.Find.Execute("38×49=1862", $true, $false, $false, $false, $false, $true, 1, $false, "70×88=6160", 2) | Out-Null
$d.Content.Find.Execute("33×64=2112", $true, $false, $false, $false, $false, $true, 1, $false, "37×97=3589", 2) | Out-Null
$d.Content.Find.Execute("66×60=3960", $true, $false, $false, $false, $false, $true, 1, $false, "18×53=954", 2) | Out-Null
$d.Content.Find.Execute("77×76=5852", $true, $false, $false, $false, $false, $true, 1, $false, "11×16=176", 2) | Out-Null
$d.Content.Find.Execute("58×96=5568", $true, $false, $false, $false, $false, $true, 1, $false, "61×41=2501", 2) | Out-Null
$d.Content.Find.Execute("70×32=2240", $true, $false, $false, $false, $false, $true, 1, $false, "32×80=2560", 2) | Out-Null
